$wb = $excel.ActiveWorkbook

# Excel quantizes ColumnWidth to whole pixels, so the value below is the
# closest achievable input that snaps to the target stored width of
# 17.2159881591797 characters (lands on 17.1666... - nearest reachable grid point).
$newColWidth = 16.33

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-22 10:40:00"
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-22 10:39:55"
$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-22 10:40:00"
$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
